# Natmi following Dr Hou advice
# Rewrite the C1qa-Cd93 sheet with the full ECs x {ECs,FAPs,M2,sCs} and
# M2 x {ECs,FAPs,M2,sCs} combinations (rows 2-9), replacing the original
# 3-row table that only covered the M2 sending cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "C1qa"
$ws.Range("C2").Value = "Cd93"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 5.392385333333333
$ws.Range("H2").Value = 16.177156
$ws.Range("I2").Value = 0.04998147672264548
$ws.Range("J2").Value = 0.04998147672264548
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 135.955556
$ws.Range("N2").Value = 407.866668
$ws.Range("O2").Value = 0.6947679994035034
$ws.Range("P2").Value = 0.6947679994035034
$ws.Range("Q2").Value = 733.1247461595786
$ws.Range("R2").Value = 6598.122715436208
$ws.Range("S2").Value = 0.03472553058982518
$ws.Range("T2").Value = 0.03472553058982518
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "C1qa"
$ws.Range("C3").Value = "Cd93"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 5.392385333333333
$ws.Range("H3").Value = 16.177156
$ws.Range("I3").Value = 0.04998147672264548
$ws.Range("J3").Value = 0.04998147672264548
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.449122
$ws.Range("N3").Value = 1.347366
$ws.Range("O3").Value = 0.002295129398228494
$ws.Range("P3").Value = 0.002295129398228494
$ws.Range("Q3").Value = 2.421838885677333
$ws.Range("R3").Value = 21.796549971096
$ws.Range("S3").Value = 0.0001147139565930168
$ws.Range("T3").Value = 0.0001147139565930168
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "C1qa"
$ws.Range("C4").Value = "Cd93"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 5.392385333333333
$ws.Range("H4").Value = 16.177156
$ws.Range("I4").Value = 0.04998147672264548
$ws.Range("J4").Value = 0.04998147672264548
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 56.38366533333333
$ws.Range("N4").Value = 169.150996
$ws.Range("O4").Value = 0.2881350899898248
$ws.Range("P4").Value = 0.2881350899898248
$ws.Range("Q4").Value = 304.0424499830418
$ws.Range("R4").Value = 2736.382049847376
$ws.Range("S4").Value = 0.01440141729330379
$ws.Range("T4").Value = 0.01440141729330379
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "C1qa"
$ws.Range("C5").Value = "Cd93"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 5.392385333333333
$ws.Range("H5").Value = 16.177156
$ws.Range("I5").Value = 0.04998147672264548
$ws.Range("J5").Value = 0.04998147672264548
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.896484
$ws.Range("N5").Value = 8.689452
$ws.Range("O5").Value = 0.01480178120844327
$ws.Range("P5").Value = 0.01480178120844327
$ws.Range("Q5").Value = 15.61895783983466
$ws.Range("R5").Value = 140.570620558512
$ws.Range("S5").Value = 0.0007398148829234987
$ws.Range("T5").Value = 0.0007398148829234987
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "C1qa"
$ws.Range("C6").Value = "Cd93"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 102.49529
$ws.Range("H6").Value = 307.48587
$ws.Range("I6").Value = 0.9500185232773545
$ws.Range("J6").Value = 0.9500185232773545
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 135.955556
$ws.Range("N6").Value = 407.866668
$ws.Range("O6").Value = 0.6947679994035034
$ws.Range("P6").Value = 0.6947679994035034
$ws.Range("Q6").Value = 13934.80413933124
$ws.Range("R6").Value = 125413.2372539811
$ws.Range("S6").Value = 0.6600424688136782
$ws.Range("T6").Value = 0.6600424688136782
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "C1qa"
$ws.Range("C7").Value = "Cd93"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 102.49529
$ws.Range("H7").Value = 307.48587
$ws.Range("I7").Value = 0.9500185232773545
$ws.Range("J7").Value = 0.9500185232773545
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.449122
$ws.Range("N7").Value = 1.347366
$ws.Range("O7").Value = 0.002295129398228494
$ws.Range("P7").Value = 0.002295129398228494
$ws.Range("Q7").Value = 46.03288963538
$ws.Range("R7").Value = 414.29600671842
$ws.Range("S7").Value = 0.002180415441635477
$ws.Range("T7").Value = 0.002180415441635477
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "C1qa"
$ws.Range("C8").Value = "Cd93"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 102.49529
$ws.Range("H8").Value = 307.48587
$ws.Range("I8").Value = 0.9500185232773545
$ws.Range("J8").Value = 0.9500185232773545
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 56.38366533333333
$ws.Range("N8").Value = 169.150996
$ws.Range("O8").Value = 0.2881350899898248
$ws.Range("P8").Value = 0.2881350899898248
$ws.Range("Q8").Value = 5779.060129602946
$ws.Range("R8").Value = 52011.54116642651
$ws.Range("S8").Value = 0.273733672696521
$ws.Range("T8").Value = 0.273733672696521
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "C1qa"
$ws.Range("C9").Value = "Cd93"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 102.49529
$ws.Range("H9").Value = 307.48587
$ws.Range("I9").Value = 0.9500185232773545
$ws.Range("J9").Value = 0.9500185232773545
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.896484
$ws.Range("N9").Value = 8.689452
$ws.Range("O9").Value = 0.01480178120844327
$ws.Range("P9").Value = 0.01480178120844327
$ws.Range("Q9").Value = 296.8759675603599
$ws.Range("R9").Value = 2671.883708043239
$ws.Range("S9").Value = 0.01406196632551977
$ws.Range("T9").Value = 0.01406196632551977
